$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column E entries (new shared strings are appended in this order:
# Ray, ?Casey, ?#Denvendra, ?Shirley)
$ws.Range("E13").Value = "Ray"
$ws.Range("E12").Value = "?Casey"
$ws.Range("E11").Value = "?#Denvendra"
$ws.Range("E10").Value = "?Shirley"
$ws.Range("E21").Value = "Arbinnav"

# Remove the now-obsolete "Shirley Tsang" entry in D26
$ws.Range("D26").ClearContents()

# Update selected cell to match the author's saved selection
$ws.Range("D14").Select()
